# Invalid Login Test and Reading From Excel
# Add a second worksheet "InvalidLogin" after the existing "ValidLogin" sheet,
# populate it with header/credential data, make it the active/selected sheet,
# and bump its zoom level.

$wb = $excel.ActiveWorkbook

# Existing sheet stays first; new sheet gets placed right after it.
$validLogin = $wb.Worksheets.Item(1)

$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "abc"
$invalidLogin.Range("B2").Value = "xyz"

# Match the source selection state (cursor parked on A3) and become the
# active/visible tab.
[void]$invalidLogin.Range("A3").Select()
$excel.ActiveWindow.Zoom = 220
